$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Product Backlog sheet: remove "Have realistic options..." row (5),
#    add new "Create options menu" row (7)
# ---------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Product Backlog")
$backlog.Range("A5:B5").ClearContents()
$backlog.Range("A7").Value = "Create options menu"
$backlog.Range("B7").Value = 4
$backlog.Activate()
$backlog.Range("A5:B5").Select()

# ---------------------------------------------------------------------
# 2. Add a new "Sprint 4" sheet right after "Product Backlog"
# ---------------------------------------------------------------------
$sprint4 = $wb.Worksheets.Add($null, $backlog)
$sprint4.Name = "Sprint 4"

# Column widths roughly matching the other sprint sheets
$sprint4.Columns.Item(1).ColumnWidth = 44.14
$sprint4.Range("B1:O1").ColumnWidth = 9.57

# Row 1: date headers (style copied from Sprint 1's date header row)
$dates = @(44121,44122,44123,44124,44125,44126,44127,44128,44129,44130,44131,44132,44133,44134)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $col = 2 + $i
    $sprint4.Cells.Item(1, $col).Value = $dates[$i]
    $sprint4.Cells.Item(1, $col).NumberFormat = "m/d/yyyy"
}

# Row 2: Optimal Trend
$sprint4.Range("A2").Value = "Optimal Trend"
$sprint4.Range("B2").Value = 14
$sprint4.Range("C2").Formula = "=B2 - `$B`$5"
$sprint4.Range("D2").Formula = "=C2 - `$B`$5"
$sprint4.Range("E2:N2").Formula = "=D2 - `$B`$5"
$sprint4.Range("O2").Value = 0
$sprint4.Range("C2:O2").NumberFormat = "0.00"

# Row 3: Actual Trend
$sprint4.Range("A3").Value = "Actual Trend "
$sprint4.Range("B3").Value = 14
$sprint4.Range("C3:G3").Value = 14
$sprint4.Range("H3").Value = 10
$sprint4.Range("I3").Value = 10
$sprint4.Range("J3").Value = 10
$sprint4.Range("K3").Value = 10
$sprint4.Range("L3").Value = 9
$sprint4.Range("M3").Value = 9
$sprint4.Range("N3").Value = 8
$sprint4.Range("O3").Value = 6

# Row 5: Hours Per Day
$sprint4.Range("A5").Value = "Hours Per Day"
$sprint4.Range("B5").Formula = "=14/13"

# Row 7-10: Use case summary table
$sprint4.Range("A7").Value = "Use Cases"
$sprint4.Range("B7").Value = "Estimates"
$sprint4.Range("C7").Value = "Status"

$sprint4.Range("A8").Value = "Create news lines"
$sprint4.Range("B8").Value = 4
$sprint4.Range("C8").Value = "Done"

$sprint4.Range("A9").Value = "Have realistic options in the store to slow the spread"
$sprint4.Range("B9").Value = 5
$sprint4.Range("C9").Value = "In progress"

$sprint4.Range("A10").Value = "Gain money as you play "
$sprint4.Range("B10").Value = 5
$sprint4.Range("C10").Value = "In progress"

# Copy the "Done" formatting (green fill) from an existing sheet
$otherSprint = $wb.Worksheets.Item("Sprint 1")
$otherSprint.Range("C9").Copy()
$sprint4.Range("C8").PasteSpecial(-4122)

# New yellow "In progress" fill
$sprint4.Range("C9:C10").Interior.Color = 65535
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Add a burndown chart to the Sprint 4 sheet
# ---------------------------------------------------------------------
$chartObj = $sprint4.ChartObjects().Add(210, 50, 390, 231)
$chartObj.Name = "Chart 1"
$chartObj.Chart.ChartType = 4
$chartObj.Chart.SetSourceData($sprint4.Range("A2:O3"))
$s1 = $chartObj.Chart.SeriesCollection(1)
$s1.XValues = $sprint4.Range("B1:O1")
$s2 = $chartObj.Chart.SeriesCollection(2)
$s2.XValues = $sprint4.Range("B1:O1")
$chartObj.Chart.HasTitle = $true
$chartObj.Chart.ChartTitle.Text = "Sprint 4 Burndown"

# ---------------------------------------------------------------------
# 4. Selection tweaks on the other sheets
# ---------------------------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint 1")
$sprint1.Activate()
$sprint1.Range("B5").Select()

$sprint3 = $wb.Worksheets.Item("Sprint 3")
$sprint3.Activate()
$sprint3.Range("N12").Select()

# ---------------------------------------------------------------------
# 5. Make "Sprint 4" the active sheet / selection, matching the target
# ---------------------------------------------------------------------
$sprint4.Activate()
$sprint4.Range("A10").Select()
